$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.782.40"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "1.656.25"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.20"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3828"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3615"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.08"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.258"
$ws.Range("E10").Value = "  +3.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08225"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.70"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.550"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.474"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001241"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "1.635.16"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.84"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06979"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.783"
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.77"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.575"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.782.49"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.089"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.33"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.99"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.229"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.52"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "1.820.18"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.960"
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.175"
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.87"
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02834"
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.170"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2522"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08836"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07172"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.17"
$ws.Range("E41").Value = "  +9.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7075"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.343"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.05"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6550"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.333"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.962"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07975"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.65"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("E51").Value = "  +0.73%  "
